$d = $word.ActiveDocument

$replacements = @(
    @{old = "787×5="; new = "810×7="},
    @{old = "921×7="; new = "659×9="},
    @{old = "469×3="; new = "909×8="},
    @{old = "597×4="; new = "677×4="},
    @{old = "520×8="; new = "389×2="},
    @{old = "648×6="; new = "188×9="},
    @{old = "825×6="; new = "907×3="},
    @{old = "450×9="; new = "124×9="},
    @{old = "272×8="; new = "174×6="},
    @{old = "171×9="; new = "844×9="},
    @{old = "900×7="; new = "578×7="},
    @{old = "855×5="; new = "525×3="},
    @{old = "415×7="; new = "787×4="},
    @{old = "823×6="; new = "503×5="},
    @{old = "389×8="; new = "604×7="},
    @{old = "104×4="; new = "895×3="},
    @{old = "866×2="; new = "423×4="},
    @{old = "670×2="; new = "275×8="},
    @{old = "903×4="; new = "910×3="},
    @{old = "905×8="; new = "677×8="},
    @{old = "760×4="; new = "133×6="},
    @{old = "671×4="; new = "616×6="},
    @{old = "750×7="; new = "360×5="},
    @{old = "507×3="; new = "440×2="},
    @{old = "503×3="; new = "704×7="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.new, 2)
}
